$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the two new worksheets ("Client" and "Nominee") at the end of the
#    workbook (after the existing "Agents" sheet).
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$client = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$client.Name = "Client"

$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$nominee = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet2)
$nominee.Name = "Nominee"

# ---------------------------------------------------------------------------
# 2. Populate "Client" sheet header row + one data row.
#    (Filled column-by-column, header then data, to match the shared-string
#    insertion order of the authored workbook. The "Choose file" image path
#    value is deliberately left for later - it is set last, below.)
# ---------------------------------------------------------------------------
$client.Range("A1").Value = "Client Password"
$client.Range("B1").Value = "Name"
$client.Range("B2").Value = "Yogi"

$client.Range("C1").Value = "Choose file"

$client.Range("D1").Value = "Gender"
$client.Range("D2").Value = "Male"

$client.Range("E1").Value = "Birth date "
$client.Range("E2").Value = 35068
$client.Range("E2").NumberFormat = "mm-dd-yy"

$client.Range("F1").Value = "Marital status"
$client.Range("F2").Value = "Unmarried"

$client.Range("G1").Value = "National ID"
$client.Range("G2").Value = 1234

$client.Range("H1").Value = "Phone"
$client.Range("H2").Value = 7894561238

$client.Range("I1").Value = "Address"
$client.Range("I2").Value = "Bangalore"

$client.Range("J1").Value = "Policy ID"
$client.Range("J2").Value = "Bangalore"

$client.Range("K1").Value = "Agent ID"
$client.Range("K2").Value = 555

$client.Range("A2").Value = 777

# ---------------------------------------------------------------------------
# 3. Populate "Nominee" sheet header row + one data row.
# ---------------------------------------------------------------------------
$nominee.Range("A1").Value = "Name"
$nominee.Range("A2").Value = "Rashmi"

$nominee.Range("B1").Value = "Gender"
$nominee.Range("B2").Value = "Female"

$nominee.Range("C1").Value = "Birth date"
$nominee.Range("C2").Value = 35547

$nominee.Range("D1").Value = "National ID"
$nominee.Range("D2").Value = 345

$nominee.Range("E1").Value = "Relationship"
$nominee.Range("E2").Value = "Friend"

$nominee.Range("F1").Value = "Priority"
$nominee.Range("F2").Value = 2

$nominee.Range("G1").Value = "Phone"
$nominee.Range("G2").Value = 7894568564

# Re-use the exact same style object for both date cells (keeps a single new
# cellXfs entry instead of two near-duplicates).
$client.Range("E2").Copy() | Out-Null
$nominee.Range("C2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 4. Set the "Choose file" image path last (matches shared-string order).
# ---------------------------------------------------------------------------
$client.Range("C2").Value = "C:\Users\sunit\Downloads\1.png"

# ---------------------------------------------------------------------------
# 5. View state: selections + active sheet/tab.
# ---------------------------------------------------------------------------
$addClient = $wb.Worksheets.Item("AddClient")
$addClient.Activate() | Out-Null
$addClient.Range("A4").Select() | Out-Null

$nominee.Activate() | Out-Null
$nominee.Range("D10").Select() | Out-Null

$client.Activate() | Out-Null
$client.Range("L2").Select() | Out-Null
